$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header C1: audioFalse -> currentPhase
$ws.Range("C1").Value = "currentPhase"

# C2 and C3: specific training audio files -> train1P2
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
